$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 1886.5518
$ws.Range("I62").Value = 1744.2858
$ws.Range("J62").Value = 1931.8182
$ws.Range("K62").Value = 1744.2858
$ws.Range("L62").Value = 1931.8182
$ws.Range("M62").Value = -1120.2858
$ws.Range("N62").Value = -3179.8182

# Row 65
$ws.Range("H65").Value = 1886.5518
$ws.Range("I65").Value = 1744.2858
$ws.Range("J65").Value = 1931.8182
$ws.Range("K65").Value = 8721.429
$ws.Range("L65").Value = 9659.091
$ws.Range("M65").Value = -5601.429
$ws.Range("N65").Value = -15899.091

# Row 86
$ws.Range("H86").Value = 7799.364
$ws.Range("I86").Value = 1896.5
$ws.Range("J86").Value = 9111.111000000001
$ws.Range("K86").Value = 1896.5
$ws.Range("L86").Value = 9111.111000000001
$ws.Range("M86").Value = -773.5
$ws.Range("N86").Value = -11357.111

# Row 89
$ws.Range("H89").Value = 7799.364
$ws.Range("I89").Value = 1896.5
$ws.Range("J89").Value = 9111.111000000001
$ws.Range("K89").Value = 9482.5
$ws.Range("L89").Value = 45555.55500000001
$ws.Range("M89").Value = -3866.5
$ws.Range("N89").Value = -56787.55500000001

# Row 105
$ws.Range("H105").Value = 48000
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 48000
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 48000
$ws.Range("N105").Value = -54988

# Row 116
$ws.Range("H116").Value = 11416.5
$ws.Range("I116").Value = 3876.25
$ws.Range("J116").Value = 14432.6
$ws.Range("K116").Value = 3876.25
$ws.Range("L116").Value = 14432.6
$ws.Range("M116").Value = -434.25
$ws.Range("N116").Value = -21316.6

# Row 118
$ws.Range("H118").Value = 726.25
$ws.Range("I118").Value = 415.7143
$ws.Range("J118").Value = 2900
$ws.Range("K118").Value = 1247.1429
$ws.Range("L118").Value = 8700
$ws.Range("M118").Value = 409.8571000000002
$ws.Range("N118").Value = -12014

# Row 129
$ws.Range("H129").Value = 766.1622
$ws.Range("I129").Value = 452.25
$ws.Range("J129").Value = 916.84
$ws.Range("K129").Value = 1356.75
$ws.Range("L129").Value = 2750.52
$ws.Range("M129").Value = 3643.25
$ws.Range("N129").Value = -12750.52

# Row 135
$ws.Range("H135").Value = 330.28
$ws.Range("I135").Value = 330.28
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 2972.52
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -437.5199999999995

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 290630.44
$ws.Range("I32").Value = 322436.28
$ws.Range("J32").Value = 12329.25
$ws.Range("K32").Value = 322436.28
$ws.Range("L32").Value = 12329.25
$ws.Range("M32").Value = -322149.28

# Row 44
$ws.Range("H44").Value = 18761
$ws.Range("I44").Value = 44
$ws.Range("J44").Value = 25000
$ws.Range("K44").Value = 44
$ws.Range("L44").Value = 25000
$ws.Range("M44").Value = 444
$ws.Range("N44").Value = -25976

# Row 45
$ws.Range("H45").Value = 1984.3572
$ws.Range("I45").Value = 1618.5714
$ws.Range("J45").Value = 2350.1428
$ws.Range("K45").Value = 1618.5714
$ws.Range("L45").Value = 2350.1428
$ws.Range("M45").Value = -1241.5714
$ws.Range("N45").Value = -3104.1428

# Row 61
$ws.Range("H61").Value = 193527.75
$ws.Range("I61").Value = 1180.7609
$ws.Range("J61").Value = 1668188
$ws.Range("K61").Value = 1180.7609
$ws.Range("L61").Value = 1668188
$ws.Range("M61").Value = -968.7609
$ws.Range("N61").Value = -1668612

# Row 74
$ws.Range("H74").Value = 6634.737
$ws.Range("I74").Value = 1189.0834
$ws.Range("J74").Value = 15970.143
$ws.Range("K74").Value = 1189.0834
$ws.Range("L74").Value = 15970.143
$ws.Range("M74").Value = -315.0834

# Row 77
$ws.Range("H77").Value = 6634.737
$ws.Range("I77").Value = 1189.0834
$ws.Range("J77").Value = 15970.143
$ws.Range("K77").Value = 5945.416999999999
$ws.Range("L77").Value = 79850.715
$ws.Range("M77").Value = -1577.416999999999

# Row 136
$ws.Range("H136").Value = 193527.75
$ws.Range("I136").Value = 1180.7609
$ws.Range("J136").Value = 1668188
$ws.Range("K136").Value = 3542.2827
$ws.Range("L136").Value = 5004564
$ws.Range("M136").Value = -992.2826999999997
$ws.Range("N136").Value = -5009664

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1995.0834
$ws.Range("I86").Value = 1166.909
$ws.Range("J86").Value = 3296.5
$ws.Range("K86").Value = 1166.909
$ws.Range("L86").Value = 3296.5
$ws.Range("M86").Value = -43.90900000000011
$ws.Range("N86").Value = -5542.5

# Row 89
$ws.Range("H89").Value = 1995.0834
$ws.Range("I89").Value = 1166.909
$ws.Range("J89").Value = 3296.5
$ws.Range("K89").Value = 5834.545
$ws.Range("L89").Value = 16482.5
$ws.Range("M89").Value = -218.5450000000001
$ws.Range("N89").Value = -27714.5

# Row 92
$ws.Range("H92").Value = 38666.668
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 38666.668
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 38666.668
$ws.Range("N92").Value = -43658.668

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 1592.25
$ws.Range("I68").Value = 763.25
$ws.Range("J68").Value = 3250.25
$ws.Range("K68").Value = 2289.75
$ws.Range("L68").Value = 9750.75
$ws.Range("M68").Value = -1478.75
$ws.Range("N68").Value = -11372.75

# Row 69
$ws.Range("H69").Value = 286.66666
$ws.Range("I69").Value = 286.66666
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 859.9999799999999
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -48.99997999999994
$ws.Range("N69").ClearContents()

# Row 71
$ws.Range("H71").Value = 1592.25
$ws.Range("I71").Value = 763.25
$ws.Range("J71").Value = 3250.25
$ws.Range("K71").Value = 6869.25
$ws.Range("L71").Value = 29252.25
$ws.Range("M71").Value = -2813.25
$ws.Range("N71").Value = -37364.25

# Row 72
$ws.Range("H72").Value = 286.66666
$ws.Range("I72").Value = 286.66666
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 2579.99994
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = 1476.00006
$ws.Range("N72").ClearContents()

# Row 107
$ws.Range("H107").Value = 750.1111
$ws.Range("I107").Value = 131.66667
$ws.Range("J107").Value = 1059.3334
$ws.Range("K107").Value = 395.00001
$ws.Range("L107").Value = 3178.0002
$ws.Range("M107").Value = 1524.99999
$ws.Range("N107").Value = -7018.0002

# Row 113
$ws.Range("H113").Value = 656830.2
$ws.Range("I113").Value = 458.86667
$ws.Range("J113").Value = 1125666.9
$ws.Range("K113").Value = 1376.60001
$ws.Range("L113").Value = 3377000.7
$ws.Range("M113").Value = 793.3999899999999
$ws.Range("N113").Value = -3381340.7

# Row 122
$ws.Range("H122").Value = 619.1622
$ws.Range("I122").Value = 234.89473
$ws.Range("J122").Value = 1024.7778
$ws.Range("K122").Value = 2114.05257
$ws.Range("L122").Value = 9223.0002
$ws.Range("M122").Value = 335.9474299999997
$ws.Range("N122").Value = -14123.0002

# Row 132
$ws.Range("H132").Value = 3483367.5
$ws.Range("I132").Value = 4348306.5
$ws.Range("J132").Value = 167767.5
$ws.Range("K132").Value = 39134758.5
$ws.Range("L132").Value = 1509907.5
$ws.Range("M132").Value = -39132228.5
$ws.Range("N132").Value = -1514967.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 32.24138
$ws.Range("I2").Value = 23.4
$ws.Range("J2").Value = 51.88889
$ws.Range("K2").Value = 23.4
$ws.Range("L2").Value = 51.88889
$ws.Range("M2").Value = 89.59999999999999
$ws.Range("N2").Value = -277.88889

# Row 109
$ws.Range("H109").Value = 16190
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 16190
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 16190
$ws.Range("N109").Value = -18270

# Row 132
$ws.Range("H132").Value = 68560.53
$ws.Range("I132").Value = 125591.336
$ws.Range("J132").Value = 4400.875
$ws.Range("K132").Value = 376774.008
$ws.Range("L132").Value = 13202.625
$ws.Range("M132").Value = -374244.008
$ws.Range("N132").Value = -18262.625
